$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 319.073924
$ws.Range("H2").Value = 957.221772
$ws.Range("I2").Value = 0.6828333423212949
$ws.Range("J2").Value = 0.6828333423212949
$ws.Range("M2").Value = 0.09834766666666667
$ws.Range("N2").Value = 0.295043
$ws.Range("O2").Value = 0.2818566198948398
$ws.Range("P2").Value = 0.2818566198948398
$ws.Range("Q2").Value = 31.38017591957733
$ws.Range("R2").Value = 282.421583276196
$ws.Range("S2").Value = 0.1924610978181762
$ws.Range("T2").Value = 0.1924610978181763

$ws.Range("G3").Value = 319.073924
$ws.Range("H3").Value = 957.221772
$ws.Range("I3").Value = 0.6828333423212949
$ws.Range("J3").Value = 0.6828333423212949
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.2505803333333333
$ws.Range("N3").Value = 0.751741
$ws.Range("O3").Value = 0.7181433801051602
$ws.Range("P3").Value = 0.7181433801051602
$ws.Range("Q3").Value = 79.95365023389466
$ws.Range("R3").Value = 719.582852105052
$ws.Range("S3").Value = 0.4903722445031187
$ws.Range("T3").Value = 0.4903722445031187

$ws.Range("G4").Value = 140.4344916666667
$ws.Range("H4").Value = 421.303475
$ws.Range("I4").Value = 0.3005364779415257
$ws.Range("J4").Value = 0.3005364779415257
$ws.Range("M4").Value = 0.09834766666666667
$ws.Range("N4").Value = 0.295043
$ws.Range("O4").Value = 0.2818566198948398
$ws.Range("P4").Value = 0.2818566198948398
$ws.Range("Q4").Value = 13.81140457493611
$ws.Range("R4").Value = 124.302641174425
$ws.Range("S4").Value = 0.08470819582769851
$ws.Range("T4").Value = 0.08470819582769852

$ws.Range("G5").Value = 140.4344916666667
$ws.Range("H5").Value = 421.303475
$ws.Range("I5").Value = 0.3005364779415257
$ws.Range("J5").Value = 0.3005364779415257
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.2505803333333333
$ws.Range("N5").Value = 0.751741
$ws.Range("O5").Value = 0.7181433801051602
$ws.Range("P5").Value = 0.7181433801051602
$ws.Range("Q5").Value = 35.19012173333056
$ws.Range("R5").Value = 316.711095599975
$ws.Range("S5").Value = 0.2158282821138272
$ws.Range("T5").Value = 0.2158282821138272

$ws.Range("G6").Value = 7.479044333333333
$ws.Range("H6").Value = 22.437133
$ws.Range("I6").Value = 0.01600550986892662
$ws.Range("J6").Value = 0.01600550986892662
$ws.Range("M6").Value = 0.09834766666666667
$ws.Range("N6").Value = 0.295043
$ws.Range("O6").Value = 0.2818566198948398
$ws.Range("P6").Value = 0.2818566198948398
$ws.Range("Q6").Value = 0.7355465590798889
$ws.Range("R6").Value = 6.619919031719
$ws.Range("S6").Value = 0.004511258911349156
$ws.Range("T6").Value = 0.004511258911349157

$ws.Range("G7").Value = 7.479044333333333
$ws.Range("H7").Value = 22.437133
$ws.Range("I7").Value = 0.01600550986892662
$ws.Range("J7").Value = 0.01600550986892662
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.2505803333333333
$ws.Range("N7").Value = 0.751741
$ws.Range("O7").Value = 0.7181433801051602
$ws.Range("P7").Value = 0.7181433801051602
$ws.Range("Q7").Value = 1.874101422061444
$ws.Range("R7").Value = 16.866912798553
$ws.Range("S7").Value = 0.01149425095757746
$ws.Range("T7").Value = 0.01149425095757746

$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0.3333333333333333
$ws.Range("G8").Value = 0.2918953333333333
$ws.Range("H8").Value = 0.875686
$ws.Range("I8").Value = 0.0006246698682528143
$ws.Range("J8").Value = 0.0006246698682528143
$ws.Range("M8").Value = 0.09834766666666667
$ws.Range("N8").Value = 0.295043
$ws.Range("O8").Value = 0.2818566198948398
$ws.Range("P8").Value = 0.2818566198948398
$ws.Range("Q8").Value = 0.02870722494422222
$ws.Range("R8").Value = 0.258365024498
$ws.Range("S8").Value = 0.0001760673376158931
$ws.Range("T8").Value = 0.0001760673376158932

$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0.3333333333333333
$ws.Range("G9").Value = 0.2918953333333333
$ws.Range("H9").Value = 0.875686
$ws.Range("I9").Value = 0.0006246698682528143
$ws.Range("J9").Value = 0.0006246698682528143
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.2505803333333333
$ws.Range("N9").Value = 0.751741
$ws.Range("O9").Value = 0.7181433801051602
$ws.Range("P9").Value = 0.7181433801051602
$ws.Range("Q9").Value = 0.07314322992511112
$ws.Range("R9").Value = 0.658289069326
$ws.Range("S9").Value = 0.0004486025306369212
$ws.Range("T9").Value = 0.0004486025306369212
